$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers - add new columns P and Q
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: update I, K, M columns and add P, Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I column: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K column: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M column: 1 -> 2
    $ws.Cells.Item($r, 16).Value = 2   # P column: new
    $ws.Cells.Item($r, 17).Value = 2   # Q column: new
}
